# "tuesday titrations 4 runs"
# Corrects the CRM standard value used on 2022-03-21 (row 70) and appends
# the fourth Tuesday titration run (2022-03-22, row 71).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 70: CRM value for the 3/21 run was re-read; fix C70 -----------
# D70's "=100*(B70-C70)/C70" formula is left untouched and recalculates
# automatically from the corrected CRM value.
$ws.Range("C70").Value = 2224.4699999999998

# --- Row 71: new titration run recorded 2022-03-22 ----------------------
$ws.Range("A71").Value = 20220322
$ws.Range("B71").Value = 2223.1790210232102
$ws.Range("C71").Value = 2224.4699999999998
$ws.Range("D71").Formula = "=100*(B71-C71)/C71"
$ws.Range("E71").Value = 180
$ws.Range("F71").Value = "CRM OPENED 20220318"

# --- View state left by the user's session when the workbook was saved --
$ws.Application.ActiveWindow.ScrollRow = 57
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C69").Select()
